$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J column (the k values)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary stats with labels
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the B14:B17 cells - bold, size 12 font, vertical center alignment, row height 15.6
$range = $ws.Range("B14:B17")
$font = $range.Font
$font.Bold = $true
$font.Size = 12
$range.VerticalAlignment = -4108

$ws.Range("A14:B17").RowHeight = 15.6

# Page setup: A4 paper, portrait orientation (matches Greek-locale Excel defaults)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A14:B17").Select()
